# "Generate Report for Handback"
# Two files (8b48e3c5-... and a302827e-...) have now been handed back
# ("Handed back: in sync with en-US") and move to the top of the report,
# pushing the still-in-progress files down. The zh-cn / de-de detail
# sheets gain "Latest Target File" / "Latest Handback File" data (with
# hyperlinks) for the two handed-back rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Source URLs (GitHub blob links) reused by the hyperlinks we (re)create.
# ---------------------------------------------------------------------
$md_616ea4f0 = "https://github.com/OpenLocalizationTest/oltest/blob/55b6d95e99c5d39c72e4e6cbc35cf313abf34108/e2e/616ea4f0-fd36-413d-8562-fab642df474d.md"
$md_26b0d971 = "https://github.com/OpenLocalizationTest/oltest/blob/ac5e95b52f0fd4c1a305bdeea5865fe14509b907/e2e/26b0d971-b515-4944-8c58-8d0385e557a2.md"
$md_8b48e3c5 = "https://github.com/OpenLocalizationTest/oltest/blob/97184a383b43ea9fbcd74a48483552b91b078627/e2e/8b48e3c5-d42d-43a3-8953-54065b6ba559.md"
$md_a302827e = "https://github.com/OpenLocalizationTest/oltest/blob/97184a383b43ea9fbcd74a48483552b91b078627/e2e/a302827e-b58a-4b6e-9062-16ccbdf64fbf.md"

$xlfZh_616ea4f0 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6545c00ead86c6e99b998e9795cb2bb44fe9f275/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/616ea4f0-fd36-413d-8562-fab642df474d.9c9a7465d6cee84c3ab8b5e5d1861567d240f28c.zh-cn.xlf"
$xlfZh_26b0d971 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3aeeaa6436dcc1aedaaf70d827e121e41aff341/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/26b0d971-b515-4944-8c58-8d0385e557a2.60c9781ba08a4f870ec94a08de4cbbb1d02af3c1.zh-cn.xlf"
$xlfZh_8b48e3c5 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3aeeaa6436dcc1aedaaf70d827e121e41aff341/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.zh-cn.xlf"
$xlfZh_a302827e = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3aeeaa6436dcc1aedaaf70d827e121e41aff341/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.zh-cn.xlf"

$xlfDe_616ea4f0 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70dfafc4ed3491729b8781312588926299a8d0ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/616ea4f0-fd36-413d-8562-fab642df474d.9c9a7465d6cee84c3ab8b5e5d1861567d240f28c.de-de.xlf"
$xlfDe_26b0d971 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f6791810ae946d6bb26f4c734d9add3558d1204/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/26b0d971-b515-4944-8c58-8d0385e557a2.60c9781ba08a4f870ec94a08de4cbbb1d02af3c1.de-de.xlf"
$xlfDe_8b48e3c5 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f6791810ae946d6bb26f4c734d9add3558d1204/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.de-de.xlf"
$xlfDe_a302827e = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f6791810ae946d6bb26f4c734d9add3558d1204/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.de-de.xlf"

$statusHandedBack = "Handed back: in sync with en-US"

# =======================================================================
# Sheet "Overview" - reorder rows: handed-back files move to the top.
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop the old hyperlinks on A2:A5 so we can lay new ones down cleanly.
$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.md"
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("D2").Value = "2016-13-12 06:13:15"

$wsOverview.Range("A3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md"
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack
$wsOverview.Range("D3").Value = "2016-13-12 06:13:15"

$wsOverview.Range("A4").Value = "616ea4f0-fd36-413d-8562-fab642df474d.md"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"
$wsOverview.Range("D4").Value = "2016-12-12 06:12:15"

$wsOverview.Range("A5").Value = "26b0d971-b515-4944-8c58-8d0385e557a2.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-13-12 06:13:15"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $md_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $md_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $md_616ea4f0, "", "", "616ea4f0-fd36-413d-8562-fab642df474d.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $md_26b0d971, "", "", "26b0d971-b515-4944-8c58-8d0385e557a2.md") | Out-Null

# =======================================================================
# Sheet "zh-cn" - same reorder, plus F/G "Latest Target / Handback File"
# data + hyperlinks for the two handed-back rows.
# =======================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A1").Hyperlinks.Delete()

# Row 2 : 8b48e3c5 (handed back)
$wsZh.Range("A2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("D2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-12 06:13:11"
$wsZh.Range("F2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.md"
$wsZh.Range("G2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-03-12 06:13:29"
$wsZh.Range("I2").Value = "Include"

# Row 3 : a302827e (handed back)
$wsZh.Range("A3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $statusHandedBack
$wsZh.Range("D3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-12 06:13:11"
$wsZh.Range("F3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md"
$wsZh.Range("G3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-03-12 06:13:29"
$wsZh.Range("I3").Value = "Include"

# Row 4 : 616ea4f0 (still in translation)
$wsZh.Range("A4").Value = "616ea4f0-fd36-413d-8562-fab642df474d.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "In Translation"
$wsZh.Range("D4").Value = "616ea4f0-fd36-413d-8562-fab642df474d.9c9a7465d6cee84c3ab8b5e5d1861567d240f28c.zh-cn.xlf"
$wsZh.Range("E4").Value = "2016-03-12 06:11:56"
$wsZh.Range("F4").Value = ""
$wsZh.Range("G4").Value = ""
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"

# Row 5 : 26b0d971 (ready for handoff)
$wsZh.Range("A5").Value = "26b0d971-b515-4944-8c58-8d0385e557a2.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "26b0d971-b515-4944-8c58-8d0385e557a2.60c9781ba08a4f870ec94a08de4cbbb1d02af3c1.zh-cn.xlf"
$wsZh.Range("E5").Value = "2016-03-12 06:13:11"
$wsZh.Range("F5").Value = ""
$wsZh.Range("G5").Value = ""
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $md_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $md_8b48e3c5, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $xlfZh_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $md_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $xlfZh_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $md_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $md_a302827e, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $xlfZh_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $md_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $xlfZh_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $md_616ea4f0, "", "", "616ea4f0-fd36-413d-8562-fab642df474d.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), $md_616ea4f0, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $xlfZh_616ea4f0, "", "", "616ea4f0-fd36-413d-8562-fab642df474d.9c9a7465d6cee84c3ab8b5e5d1861567d240f28c.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $md_26b0d971, "", "", "26b0d971-b515-4944-8c58-8d0385e557a2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B5"), $md_26b0d971, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), $xlfZh_26b0d971, "", "", "26b0d971-b515-4944-8c58-8d0385e557a2.60c9781ba08a4f870ec94a08de4cbbb1d02af3c1.zh-cn.xlf") | Out-Null

# =======================================================================
# Sheet "de-de" - same reorder, plus F/G "Latest Target / Handback File"
# data + hyperlinks for the two handed-back rows.
# =======================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A1").Hyperlinks.Delete()

# Row 2 : 8b48e3c5 (handed back)
$wsDe.Range("A2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("D2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-12 06:13:15"
$wsDe.Range("F2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.md"
$wsDe.Range("G2").Value = "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.de-de.xlf"
$wsDe.Range("H2").Value = "2016-03-12 06:13:34"
$wsDe.Range("I2").Value = "Include"

# Row 3 : a302827e (handed back)
$wsDe.Range("A3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $statusHandedBack
$wsDe.Range("D3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-12 06:13:15"
$wsDe.Range("F3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md"
$wsDe.Range("G3").Value = "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.de-de.xlf"
$wsDe.Range("H3").Value = "2016-03-12 06:13:34"
$wsDe.Range("I3").Value = "Include"

# Row 4 : 616ea4f0 (still in translation)
$wsDe.Range("A4").Value = "616ea4f0-fd36-413d-8562-fab642df474d.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "In Translation"
$wsDe.Range("D4").Value = "616ea4f0-fd36-413d-8562-fab642df474d.9c9a7465d6cee84c3ab8b5e5d1861567d240f28c.de-de.xlf"
$wsDe.Range("E4").Value = "2016-03-12 06:12:15"
$wsDe.Range("F4").Value = ""
$wsDe.Range("G4").Value = ""
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"

# Row 5 : 26b0d971 (ready for handoff)
$wsDe.Range("A5").Value = "26b0d971-b515-4944-8c58-8d0385e557a2.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "26b0d971-b515-4944-8c58-8d0385e557a2.60c9781ba08a4f870ec94a08de4cbbb1d02af3c1.de-de.xlf"
$wsDe.Range("E5").Value = "2016-03-12 06:13:15"
$wsDe.Range("F5").Value = ""
$wsDe.Range("G5").Value = ""
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $md_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $md_8b48e3c5, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $xlfDe_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $md_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $xlfDe_8b48e3c5, "", "", "8b48e3c5-d42d-43a3-8953-54065b6ba559.6969a453e2e4c62684c93bdca792f0aa8e2e76ff.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $md_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $md_a302827e, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $xlfDe_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $md_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $xlfDe_a302827e, "", "", "a302827e-b58a-4b6e-9062-16ccbdf64fbf.5d841f859edd579055c7dd3e87f00a382b9f1fa4.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $md_616ea4f0, "", "", "616ea4f0-fd36-413d-8562-fab642df474d.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), $md_616ea4f0, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $xlfDe_616ea4f0, "", "", "616ea4f0-fd36-413d-8562-fab642df474d.9c9a7465d6cee84c3ab8b5e5d1861567d240f28c.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $md_26b0d971, "", "", "26b0d971-b515-4944-8c58-8d0385e557a2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B5"), $md_26b0d971, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), $xlfDe_26b0d971, "", "", "26b0d971-b515-4944-8c58-8d0385e557a2.60c9781ba08a4f870ec94a08de4cbbb1d02af3c1.de-de.xlf") | Out-Null
